# This script rotates the species-observation data among rows 2, 3, 4 and 6
# of the active worksheet, while leaving row 5 and all "static" per-row
# columns (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY)
# untouched.
#
# Net effect (matches target diff):
#   row 6 data -> row 2
#   row 2 data -> row 3
#   row 3 data -> row 4
#   row 4 data -> row 6
#
# Note: this engine's Value property reads back oddly when accessed without
# an explicit indexer argument, so Value2 is used for both reading and
# writing cell contents (it is reliable in both directions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "moving" field values (A, B, E, F, G, H, Q, R, AC) for the
# four rows involved in the rotation before any writes happen.
$row2_A = $ws.Range("A2").Value2
$row2_B = $ws.Range("B2").Value2
$row2_E = $ws.Range("E2").Value2
$row2_F = $ws.Range("F2").Value2
$row2_G = $ws.Range("G2").Value2
$row2_H = $ws.Range("H2").Value2
$row2_Q = $ws.Range("Q2").Value2
$row2_R = $ws.Range("R2").Value2
$row2_AC = $ws.Range("AC2").Value2

$row3_A = $ws.Range("A3").Value2
$row3_B = $ws.Range("B3").Value2
$row3_E = $ws.Range("E3").Value2
$row3_F = $ws.Range("F3").Value2
$row3_G = $ws.Range("G3").Value2
$row3_H = $ws.Range("H3").Value2
$row3_Q = $ws.Range("Q3").Value2
$row3_R = $ws.Range("R3").Value2
$row3_AC = $ws.Range("AC3").Value2

$row4_A = $ws.Range("A4").Value2
$row4_B = $ws.Range("B4").Value2
$row4_E = $ws.Range("E4").Value2
$row4_F = $ws.Range("F4").Value2
$row4_G = $ws.Range("G4").Value2
$row4_H = $ws.Range("H4").Value2
$row4_Q = $ws.Range("Q4").Value2
$row4_R = $ws.Range("R4").Value2
$row4_AC = $ws.Range("AC4").Value2

$row6_A = $ws.Range("A6").Value2
$row6_B = $ws.Range("B6").Value2
$row6_E = $ws.Range("E6").Value2
$row6_F = $ws.Range("F6").Value2
$row6_G = $ws.Range("G6").Value2
$row6_H = $ws.Range("H6").Value2
$row6_Q = $ws.Range("Q6").Value2
$row6_R = $ws.Range("R6").Value2
$row6_AC = $ws.Range("AC6").Value2

# Row 2 <- old row 6
$ws.Range("A2").Value2 = $row6_A
$ws.Range("B2").Value2 = $row6_B
$ws.Range("E2").Value2 = $row6_E
$ws.Range("F2").Value2 = $row6_F
$ws.Range("G2").Value2 = $row6_G
$ws.Range("H2").Value2 = $row6_H
$ws.Range("Q2").Value2 = $row6_Q
$ws.Range("R2").Value2 = $row6_R
$ws.Range("AC2").Value2 = $row6_AC

# Row 3 <- old row 2
$ws.Range("A3").Value2 = $row2_A
$ws.Range("B3").Value2 = $row2_B
$ws.Range("E3").Value2 = $row2_E
$ws.Range("F3").Value2 = $row2_F
$ws.Range("G3").Value2 = $row2_G
$ws.Range("H3").Value2 = $row2_H
$ws.Range("Q3").Value2 = $row2_Q
$ws.Range("R3").Value2 = $row2_R
$ws.Range("AC3").Value2 = $row2_AC

# Row 4 <- old row 3
$ws.Range("A4").Value2 = $row3_A
$ws.Range("B4").Value2 = $row3_B
$ws.Range("E4").Value2 = $row3_E
$ws.Range("F4").Value2 = $row3_F
$ws.Range("G4").Value2 = $row3_G
$ws.Range("H4").Value2 = $row3_H
$ws.Range("Q4").Value2 = $row3_Q
$ws.Range("R4").Value2 = $row3_R
$ws.Range("AC4").Value2 = $row3_AC

# Row 6 <- old row 4
$ws.Range("A6").Value2 = $row4_A
$ws.Range("B6").Value2 = $row4_B
$ws.Range("E6").Value2 = $row4_E
$ws.Range("F6").Value2 = $row4_F
$ws.Range("G6").Value2 = $row4_G
$ws.Range("H6").Value2 = $row4_H
$ws.Range("Q6").Value2 = $row4_Q
$ws.Range("R6").Value2 = $row4_R
$ws.Range("AC6").Value2 = $row4_AC
